$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are stored as text (inline strings) in the workbook, even when
# they look numeric (e.g. "5.10", "0.166"). Excel auto-converts plain-looking
# numeric strings assigned to .Value into real numbers, which would both change
# the cell type and drop significant trailing zeros. To keep those cells as text
# (matching the source data) we briefly force a Text number format, assign the
# value, then restore the Normal style so no formatting residue is left behind.

$ws.Range('D2').Value = '67.922.95'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').Value = '2.539.40'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').Value = '2.538.43'
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.166'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.10'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.53'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.973.66'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000178'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('D17').Value = '67.800.04'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').Value = '2.524.17'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.83%  '
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '370.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.85%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E26').Value = '  -3.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.52%  '
$ws.Range('D28').Value = '2.667.65'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('D29').Value = '0.0₃0970'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('E30').Value = '  +4.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '541.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.87'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.92'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.18'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.64'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('E42').Value = '  -0.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.58'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.45%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '39.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('D46').Value = '0.0₆0290'
$ws.Range('E46').Value = '  +4.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '147.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.72'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.553'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0746'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.59%  '
